$wb = $excel.ActiveWorkbook

# --- Sheet1: update Facility ID values (HMSL7000x -> HMSL7500x) ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A3").Value  = "HMSL75001"
$ws1.Range("A4").Value  = "HMSL75001"
$ws1.Range("A5").Value  = "HMSL75002"
$ws1.Range("A6").Value  = "HMSL75002"
$ws1.Range("A7").Value  = "HMSL75003"
$ws1.Range("A8").Value  = "HMSL75003"
$ws1.Range("A9").Value  = "HMSL75004"
$ws1.Range("A10").Value = "HMSL75004"
$ws1.Range("A11").Value = "HMSL75005"
$ws1.Range("A12").Value = "HMSL75006"
$ws1.Range("A13").Value = "HMSL75007"
$ws1.Range("A14").Value = "HMSL75008"
$ws1.Range("A15").Value = "HMSL75009"

# --- Sheet1: DC_Center_Batch_ID for row 12 changes from 2 to 1 ---
$ws1.Range("B12").Value = 1

# --- Sheet1: scroll/selection state changes (frozen pane view moved) ---
$ws1.Activate()
$ws1.Range("B13").Select()
